$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert two new blank rows at rows 5-6 (pushes everything from
#    the old row 5 downward by two rows). Excel automatically
#    copies the row-above (row 4) formatting into the new rows,
#    which produces the "B5"/"B6" style-3 placeholder cells seen
#    in the target file.
# ---------------------------------------------------------------
$ws.Range("A5:A6").EntireRow.Insert()

# ---------------------------------------------------------------
# 2. Fill in the two newly created label rows (bold label style,
#    matching the other D5:D14 caption cells).
# ---------------------------------------------------------------
$ws.Range("D5").Value = "Từ khóa:"
$ws.Range("D5").Font.Bold = $true
$ws.Range("D6").Value = "Người nộp đơn:"
$ws.Range("D6").Font.Bold = $true

# ---------------------------------------------------------------
# 3. Fix the typo in the "Thời gian tiếp nhận:" label (was
#    previously "tThời gian tiếp nhận:") - this now lives at D11
#    after the 2-row shift (old D9).
# ---------------------------------------------------------------
$ws.Range("D11").Value = "Thời gian tiếp nhận:"

# ---------------------------------------------------------------
# 4. "Công khai:" (now at D14 after the shift, old D12) used to be
#    drawn with the plain font; make it bold like the rest of the
#    label column.
# ---------------------------------------------------------------
$ws.Range("D14").Font.Bold = $true

# ---------------------------------------------------------------
# 5. Update the sheet view: unfreeze, reselect the first data row
#    (now row 17) and re-freeze so the freeze pane follows the
#    header row, move the visible selection to D12, and bump the
#    zoom level from 115% to 175%.
# ---------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A17").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 175
